$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("titreMessage"); the old column B ("Messages en
# français" header + its message text) shifts right to become column C.
$ws.Columns("B").Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "titreMessage"

# New row 3: a second message entry.
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Salut"
$ws.Range("C3").Value = "Salut tout le monde !"

# Row 2: title for the pre-existing message (now in C2).
$ws.Range("B2").Value = "Hello world"

# Column widths.
$ws.Columns("B").ColumnWidth = 27.5
$ws.Columns("C").ColumnWidth = 76

# Header row shading (A1:C1) - yellow fill.
$ws.Range("A1:B1").Interior.Color = 65535

# Header for the message column also gets wrap text in addition to the fill.
$ws.Range("C1").Interior.Color = 65535
$ws.Range("C1").WrapText = $true

# Message column body wraps text too.
$ws.Range("C2:C3").WrapText = $true

# Selection, matching the saved state.
$ws.Range("B3").Select()

$wb.Save()
